# update data day 20/10/2024
# Append the day's scraped Facebook post URLs below the existing "link"
# header in column A, then leave the selection on the last entered cell
# (mirrors the author's saved selection state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
  "https://www.facebook.com/K14vn/posts/964043959089937?ref=embed_post",
  "https://www.facebook.com/K14vn/posts/963780079116325?ref=embed_post",
  "https://www.facebook.com/K14vn/posts/963721785788821?ref=embed_post",
  "https://www.facebook.com/K14vn/posts/963387679155565?ref=embed_post",
  "https://www.facebook.com/K14vn/posts/963266382501028?ref=embed_post",
  "https://www.facebook.com/K14vn/posts/963193552508311?ref=embed_post",
  "https://www.facebook.com/K14vn/posts/963055725855427?ref=embed_post",
  "https://www.facebook.com/K14vn/posts/963038832523783?ref=embed_post",
  "https://www.facebook.com/K14vn/posts/963026039191729?ref=embed_post",
  "https://www.facebook.com/K14vn/posts/962881379206195?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/988854093279371?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/988882469943200?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/989049856593128?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/987988396699274?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/988029190028528?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/987443800087067?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/988867133278067?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/988127853351995?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/988082110023236?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/986517283513052?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/986576653507115?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/986622820169165?ref=embed_post",
  "https://www.facebook.com/thongtinchinhphu/posts/986570580174389?ref=embed_post"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $urls[$i]
}

# Matches the saved selection in the edited workbook (A24, the last row
# written).
$ws.Range("A24").Select()

# Matches the narrower saved window width in the edited workbook
# (windowWidth 27945 -> 13545 twips).
$wb.Windows.Item(1).Width = 13545
